$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 460129
$ws.Range("D2").Value = 685402007
$ws.Range("C3").Value = 354
$ws.Range("D3").Value = 520577
$ws.Range("C9").Value = 1525
$ws.Range("D9").Value = 4243491
$ws.Range("C11").Value = 178912
$ws.Range("D11").Value = 443288343
$ws.Range("C12").Value = 331
$ws.Range("D12").Value = 1380601
$ws.Range("C14").Value = 93162
$ws.Range("D14").Value = 220558115
$ws.Range("C18").Value = 5009
$ws.Range("D18").Value = 8572122
$ws.Range("C21").Value = 111
$ws.Range("D21").Value = 332596
$ws.Range("C23").Value = 12045
$ws.Range("D23").Value = 25110145
$ws.Range("C25").Value = 109297
$ws.Range("D25").Value = 159899794
$ws.Range("C31").Value = 48316
$ws.Range("D31").Value = 116108254
$ws.Range("C32").Value = 81
$ws.Range("D32").Value = 310571
$ws.Range("C34").Value = 17565
$ws.Range("D34").Value = 40828266
$ws.Range("C37").Value = 2035
$ws.Range("D37").Value = 3645047
$ws.Range("C39").Value = 3087
$ws.Range("D39").Value = 6122992
$ws.Range("C40").Value = 134670
$ws.Range("D40").Value = 200458382
$ws.Range("C42").Value = 116
$ws.Range("D42").Value = 229620
$ws.Range("C46").Value = 1702
$ws.Range("D46").Value = 5194022
$ws.Range("C48").Value = 67620
$ws.Range("D48").Value = 167756838
$ws.Range("C51").Value = 14416
$ws.Range("D51").Value = 34752051
$ws.Range("C53").Value = 1739
$ws.Range("D53").Value = 2869497
$ws.Range("C56").Value = 3880
$ws.Range("D56").Value = 7873477
$ws.Range("C57").Value = 96948
$ws.Range("D57").Value = 144158120
$ws.Range("C63").Value = 627
$ws.Range("D63").Value = 1715076
$ws.Range("C65").Value = 41797
$ws.Range("D65").Value = 98378526
$ws.Range("C66").Value = 71
$ws.Range("D66").Value = 309094
$ws.Range("C68").Value = 16825
$ws.Range("D68").Value = 38248020
$ws.Range("C70").Value = 1790
$ws.Range("D70").Value = 3242260
$ws.Range("C74").Value = 2555
$ws.Range("D74").Value = 5183363
$ws.Range("C76").Value = 28844
$ws.Range("D76").Value = 44630124
$ws.Range("C80").Value = 11966
$ws.Range("D80").Value = 31940052
$ws.Range("C81").Value = 11
$ws.Range("D81").Value = 64360
$ws.Range("C82").Value = 8113
$ws.Range("D82").Value = 19944131
$ws.Range("C85").Value = 504
$ws.Range("D85").Value = 1000132
$ws.Range("C86").Value = 201948
$ws.Range("D86").Value = 305942269
$ws.Range("C90").Value = 781
$ws.Range("D90").Value = 2117774
$ws.Range("C91").Value = 24
$ws.Range("D91").Value = 53852
$ws.Range("C92").Value = 94963
$ws.Range("D92").Value = 224743764
$ws.Range("C95").Value = 44950
$ws.Range("D95").Value = 103113227
$ws.Range("C97").Value = 45
$ws.Range("D97").Value = 281534
$ws.Range("C98").Value = 7339
$ws.Range("D98").Value = 27724507
$ws.Range("C100").Value = 4769
$ws.Range("D100").Value = 9699285
$ws.Range("C102").Value = 44706
$ws.Range("D102").Value = 66498231
$ws.Range("C106").Value = 10870
$ws.Range("D106").Value = 18968334
$ws.Range("C108").Value = 10176
$ws.Range("D108").Value = 16992399
$ws.Range("C112").Value = 18661
$ws.Range("D112").Value = 38177287
$ws.Range("C115").Value = 4190
$ws.Range("D115").Value = 9157134
$ws.Range("C117").Value = 6083
$ws.Range("D117").Value = 13730950
$ws.Range("C122").Value = 201688
$ws.Range("D122").Value = 294195230
$ws.Range("C123").Value = 45
$ws.Range("D123").Value = 87813
$ws.Range("C128").Value = 1637
$ws.Range("D128").Value = 4579279
$ws.Range("C130").Value = 78979
$ws.Range("D130").Value = 188268229
$ws.Range("C131").Value = 191
$ws.Range("D131").Value = 694182
$ws.Range("C133").Value = 42413
$ws.Range("D133").Value = 97332975
$ws.Range("C136").Value = 1846
$ws.Range("D136").Value = 3780460
$ws.Range("C140").Value = 3961
$ws.Range("D140").Value = 8134997
$ws.Range("C142").Value = 806347
$ws.Range("D142").Value = 1288770553
$ws.Range("C147").Value = 2828
$ws.Range("D147").Value = 9942002
$ws.Range("C149").Value = 327305
$ws.Range("D149").Value = 792478746
$ws.Range("C150").Value = 1103
$ws.Range("D150").Value = 4852104
$ws.Range("C152").Value = 303598
$ws.Range("D152").Value = 688853702
$ws.Range("C155").Value = 3866
$ws.Range("D155").Value = 6531082
$ws.Range("C158").Value = 12070
$ws.Range("D158").Value = 25545550
$ws.Range("C161").Value = 55027
$ws.Range("D161").Value = 76563928
$ws.Range("C167").Value = 16525
$ws.Range("D167").Value = 26794546
$ws.Range("C169").Value = 4580
$ws.Range("D169").Value = 7320554
$ws.Range("C175").Value = 24813
$ws.Range("D175").Value = 37061954
$ws.Range("C179").Value = 10426
$ws.Range("D179").Value = 21450148
$ws.Range("C181").Value = 7331
$ws.Range("D181").Value = 13909730
$ws.Range("C184").Value = 438
$ws.Range("D184").Value = 812837
$ws.Range("C186").Value = 33529
$ws.Range("D186").Value = 81892829
$ws.Range("C193").Value = 124237
$ws.Range("D193").Value = 186302967
$ws.Range("C199").Value = 1183
$ws.Range("D199").Value = 3591053
$ws.Range("C201").Value = 51589
$ws.Range("D201").Value = 127351624
$ws.Range("C203").Value = 20414
$ws.Range("D203").Value = 48433760
$ws.Range("C205").Value = 1606
$ws.Range("D205").Value = 2806647
$ws.Range("C208").Value = 2861
$ws.Range("D208").Value = 5906761
$ws.Range("C210").Value = 329076
$ws.Range("D210").Value = 474706029
$ws.Range("C212").Value = 242
$ws.Range("D212").Value = 493588
$ws.Range("C217").Value = 1453
$ws.Range("D217").Value = 3804563
$ws.Range("C219").Value = 127556
$ws.Range("D219").Value = 300624876
$ws.Range("C220").Value = 195
$ws.Range("D220").Value = 688706
$ws.Range("C221").Value = 15
$ws.Range("D221").Value = 63332
$ws.Range("C222").Value = 50313
$ws.Range("D222").Value = 115088924
$ws.Range("C225").Value = 6637
$ws.Range("D225").Value = 12222419
$ws.Range("C228").Value = 8286
$ws.Range("D228").Value = 16011441
$ws.Range("C231").Value = 373372
$ws.Range("D231").Value = 528840476
$ws.Range("C233").Value = 356
$ws.Range("D233").Value = 684198
$ws.Range("C240").Value = 141969
$ws.Range("D240").Value = 333620170
$ws.Range("C241").Value = 213
$ws.Range("D241").Value = 792439
$ws.Range("C243").Value = 78681
$ws.Range("D243").Value = 178600684
$ws.Range("C246").Value = 6050
$ws.Range("D246").Value = 10073618
$ws.Range("C249").Value = 50
$ws.Range("D249").Value = 198522
$ws.Range("C250").Value = 10458
$ws.Range("D250").Value = 20394533
$ws.Range("C253").Value = 148681
$ws.Range("D253").Value = 219873282
$ws.Range("C259").Value = 1054
$ws.Range("D259").Value = 3083395
$ws.Range("C261").Value = 73822
$ws.Range("D261").Value = 177491783
$ws.Range("C264").Value = 19481
$ws.Range("D264").Value = 45015251
$ws.Range("C266").Value = 2407
$ws.Range("D266").Value = 4428091
$ws.Range("C268").Value = 4126
$ws.Range("D268").Value = 8350140
$ws.Range("C269").Value = 374375
$ws.Range("D269").Value = 547839253
$ws.Range("C271").Value = 375
$ws.Range("D271").Value = 813870
$ws.Range("C276").Value = 1600
$ws.Range("D276").Value = 5027305
$ws.Range("C278").Value = 148078
$ws.Range("D278").Value = 361530919
$ws.Range("C280").Value = 17
$ws.Range("D280").Value = 87065
$ws.Range("C281").Value = 104940
$ws.Range("D281").Value = 245217315
$ws.Range("C283").Value = 12
$ws.Range("D283").Value = 26500
$ws.Range("C284").Value = 3319
$ws.Range("D284").Value = 5772451
$ws.Range("C287").Value = 8558
$ws.Range("D287").Value = 17647194
